# TC05_CDS_Filter_Study-CIDR_Aggressive_Prostate_Cancer.xlsx
# - Fix the "Tumor" column of the SamplesTab Cypher query: it referenced the
#   collected `tumor` alias (from a WITH ... COLLECT(...) as tumor) instead of
#   the per-row `samp.sample_tumor_status` property. Replace the query text
#   in B3 with the corrected version.
# - Update the active selection left behind on the sheet (B3 instead of D3).
# - The corrected/longer query text wraps to more lines, so the row grows
#   taller (189 -> 204.75).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTumorQuery = "`nMATCH (s:study)<--(p:participant)<--(samp:sample)`n" +
    "WHERE s.study_name in [`"CIDR: The Genetic Basis of Aggressive Prostate Cancer: The Role of Rare Variation`"]`n" +
    "WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`n" +
    "RETURN  `n" +
    " coalesce(samp.sample_id, '') as ``Sample ID``,`n" +
    " coalesce(p.participant_id,'') as ``Participant ID``,`n" +
    " coalesce(s.study_name, '') as ``Study Name``,`n" +
    " coalesce(s.phs_accession,'') as ``Accession``,`n" +
    " coalesce(samp.sample_tumor_status,'') as ``Tumor``,`n" +
    "coalesce(samp.sample_type,'') as ``Analyte Type```n" +
    "  ORDER By samp.sample_id LIMIT 100"

$ws.Range("B3").Value = $newTumorQuery

$ws.Rows.Item(3).RowHeight = 204.75

$ws.Range("B3").Select()
